# Fruta / hortaliza, semanal
# Insert the latest week's two data points (Primera / Segunda quality) for
# Albahaca at the top of the data block (row 305), pushing the existing
# rows 305-318 down to 307-320.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 305 (existing data,
# previously rows 305:318, shifts down to 307:320).
$ws.Rows("305:306").Insert()

# New row 305 - "Primera" quality, most recent reporting date.
$ws.Range("A305").Value = 9
$ws.Range("B305").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C305").Value = "Metropolitana"
$ws.Range("D305").Value = 44610
$ws.Range("E305").Value = 13
$ws.Range("F305").Value = 100112052
$ws.Range("G305").Value = "Albahaca"
$ws.Range("H305").Value = "Sin especificar"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 160
$ws.Range("K305").Value = 3000
$ws.Range("L305").Value = 3500
$ws.Range("M305").Value = 3250
$ws.Range("N305").Value = '$/docena de matas'
$ws.Range("O305").Value = "Región Metropolitana"
$ws.Range("P305").Value = 542
$ws.Range("Q305").Value = 6
$ws.Range("R305").Value = "Hortaliza"

# New row 306 - "Segunda" quality, same reporting date.
$ws.Range("A306").Value = 9
$ws.Range("B306").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C306").Value = "Metropolitana"
$ws.Range("D306").Value = 44610
$ws.Range("E306").Value = 13
$ws.Range("F306").Value = 100112052
$ws.Range("G306").Value = "Albahaca"
$ws.Range("H306").Value = "Sin especificar"
$ws.Range("I306").Value = "Segunda"
$ws.Range("J306").Value = 52
$ws.Range("K306").Value = 2500
$ws.Range("L306").Value = 2500
$ws.Range("M306").Value = 2500
$ws.Range("N306").Value = '$/docena de matas'
$ws.Range("O306").Value = "Región Metropolitana"
$ws.Range("P306").Value = 417
$ws.Range("Q306").Value = 6
$ws.Range("R306").Value = "Hortaliza"
